$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Range('A30').Value = '2024-10-12 15:36:13'
$ws.Range('C30').Value = 0
$ws.Range('D30').Value = 0
$ws.Range('E30').Value = 0
$ws.Range('F30').Value = 0
$ws.Range('G30').Value = 0
$ws.Range('H30').Value = 0
$ws.Range('N30').Value = 10
$ws.Range('O30').Value = 10
$ws.Range('P30').Value = 1
$ws.Range('R30').Value = 5
$ws.Range('T30').Value = 30
$ws.Range('U30').NumberFormat = "@"
$ws.Range('U30').Value = '0'
$ws.Range('V30').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Crupier.xlsx'
$ws.Range('X30').Value = 'No es Simulación'
$ws.Range('Y30').Value = 0

# Row 31
$ws.Range('A31').Value = '2024-10-12 15:48:24'
$ws.Range('C31').Value = 4
$ws.Range('D31').Value = 0
$ws.Range('E31').Value = 4
$ws.Range('F31').Value = 0
$ws.Range('G31').Value = 0
$ws.Range('H31').Value = 0
$ws.Range('N31').Value = 10
$ws.Range('O31').Value = 10
$ws.Range('P31').Value = 1
$ws.Range('R31').Value = 5
$ws.Range('T31').Value = 30
$ws.Range('U31').NumberFormat = "@"
$ws.Range('U31').Value = '18%'
$ws.Range('V31').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X31').Value = 'No es Simulación'
$ws.Range('Y31').Value = 22

# Row 32
$ws.Range('A32').Value = '2024-10-12 16:22:36'
$ws.Range('C32').Value = 9
$ws.Range('D32').Value = 2
$ws.Range('E32').Value = 4
$ws.Range('F32').Value = 3
$ws.Range('G32').Value = 0
$ws.Range('H32').Value = 0
$ws.Range('N32').Value = 10
$ws.Range('O32').Value = 10
$ws.Range('P32').Value = 2
$ws.Range('R32').Value = 5
$ws.Range('T32').Value = 30
$ws.Range('U32').NumberFormat = "@"
$ws.Range('U32').Value = '32%'
$ws.Range('V32').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X32').Value = 'No es Simulación'
$ws.Range('Y32').Value = 28

# Row 33
$ws.Range('A33').Value = '2024-10-12 16:24:27'
$ws.Range('C33').Value = 2
$ws.Range('D33').Value = 0
$ws.Range('E33').Value = 2
$ws.Range('F33').Value = 0
$ws.Range('G33').Value = 0
$ws.Range('H33').Value = 0
$ws.Range('N33').Value = 10
$ws.Range('O33').Value = 10
$ws.Range('P33').Value = 1
$ws.Range('R33').Value = 5
$ws.Range('T33').Value = 30
$ws.Range('U33').NumberFormat = "@"
$ws.Range('U33').Value = '18%'
$ws.Range('V33').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X33').Value = 'No es Simulación'
$ws.Range('Y33').Value = 11

# Row 34
$ws.Range('A34').Value = '2024-10-12 16:26:14'
$ws.Range('C34').Value = 5
$ws.Range('D34').Value = 3
$ws.Range('E34').Value = 2
$ws.Range('F34').Value = 0
$ws.Range('G34').Value = 0
$ws.Range('H34').Value = 0
$ws.Range('N34').Value = 10
$ws.Range('O34').Value = 10
$ws.Range('P34').Value = 1
$ws.Range('R34').Value = 5
$ws.Range('T34').Value = 20
$ws.Range('U34').NumberFormat = "@"
$ws.Range('U34').Value = '24%'
$ws.Range('V34').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X34').Value = 'No es Simulación'
$ws.Range('Y34').Value = 21

# Row 35
$ws.Range('A35').Value = '2024-10-12 20:01:28'
$ws.Range('C35').Value = 0
$ws.Range('D35').Value = 0
$ws.Range('E35').Value = 0
$ws.Range('F35').Value = 0
$ws.Range('G35').Value = 0
$ws.Range('H35').Value = 0
$ws.Range('N35').Value = 10
$ws.Range('O35').Value = 10
$ws.Range('P35').Value = 1
$ws.Range('R35').Value = 5
$ws.Range('T35').Value = 30
$ws.Range('U35').NumberFormat = "@"
$ws.Range('U35').Value = '0'
$ws.Range('V35').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Crupier.xlsx'
$ws.Range('X35').Value = 'No es Simulación'
$ws.Range('Y35').Value = 0

# Row 36
$ws.Range('A36').Value = '2024-10-12 20:18:24'
$ws.Range('C36').Value = 10
$ws.Range('D36').Value = 2
$ws.Range('E36').Value = 8
$ws.Range('F36').Value = 0
$ws.Range('G36').Value = 0
$ws.Range('H36').Value = 0
$ws.Range('N36').Value = 10
$ws.Range('O36').Value = 10
$ws.Range('P36').Value = 1
$ws.Range('R36').Value = 5
$ws.Range('T36').Value = 30
$ws.Range('U36').NumberFormat = "@"
$ws.Range('U36').Value = '32%'
$ws.Range('V36').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X36').Value = 'No es Simulación'
$ws.Range('Y36').Value = 31

# Row 37
$ws.Range('A37').Value = '2024-10-12 20:23:07'
$ws.Range('C37').Value = 1
$ws.Range('D37').Value = 0
$ws.Range('E37').Value = 1
$ws.Range('F37').Value = 0
$ws.Range('G37').Value = 0
$ws.Range('H37').Value = 0
$ws.Range('N37').Value = 10
$ws.Range('O37').Value = 10
$ws.Range('P37').Value = 1
$ws.Range('R37').Value = 5
$ws.Range('T37').Value = 50
$ws.Range('U37').NumberFormat = "@"
$ws.Range('U37').Value = '25%'
$ws.Range('V37').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X37').Value = 'No es Simulación'
$ws.Range('Y37').Value = 4

# Row 38
$ws.Range('A38').Value = '2024-10-13 02:59:15'
$ws.Range('C38').Value = 12
$ws.Range('D38').Value = 5
$ws.Range('E38').Value = 7
$ws.Range('F38').Value = 0
$ws.Range('G38').Value = 0
$ws.Range('H38').Value = 0
$ws.Range('N38').Value = 10
$ws.Range('O38').Value = 10
$ws.Range('P38').Value = 1
$ws.Range('R38').Value = 5
$ws.Range('T38').Value = 20
$ws.Range('U38').NumberFormat = "@"
$ws.Range('U38').Value = '21%'
$ws.Range('V38').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X38').Value = 'No es Simulación'
$ws.Range('Y38').Value = 57

# Row 39
$ws.Range('A39').Value = '2024-10-13 10:23:47'
$ws.Range('C39').Value = 18
$ws.Range('D39').Value = 8
$ws.Range('E39').Value = 10
$ws.Range('F39').Value = 0
$ws.Range('G39').Value = 0
$ws.Range('H39').Value = 0
$ws.Range('N39').Value = 10
$ws.Range('O39').Value = 10
$ws.Range('P39').Value = 1
$ws.Range('R39').Value = 5
$ws.Range('T39').Value = 20
$ws.Range('U39').NumberFormat = "@"
$ws.Range('U39').Value = '41%'
$ws.Range('V39').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X39').Value = 'No es Simulación'
$ws.Range('Y39').Value = 44

# Row 40
$ws.Range('A40').Value = '2024-10-14 19:02:31'
$ws.Range('C40').Value = 0
$ws.Range('D40').Value = 0
$ws.Range('E40').Value = 0
$ws.Range('F40').Value = 0
$ws.Range('G40').Value = 0
$ws.Range('H40').Value = 0
$ws.Range('N40').Value = 10
$ws.Range('O40').Value = 10
$ws.Range('P40').Value = 1
$ws.Range('R40').Value = 5
$ws.Range('T40').Value = 50
$ws.Range('U40').NumberFormat = "@"
$ws.Range('U40').Value = '0'
$ws.Range('V40').Value = 'C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx'
$ws.Range('X40').Value = 'No es Simulación'
$ws.Range('Y40').Value = 0
